# Problem Solving #3 Step 4,5 - complete steps 4 & 5
$d = $word.ActiveDocument

# Anchor everything off the unique "My sub-goal..." paragraph in the
# 3rd "Problem Solving" write-up (finger counting problem) so we never
# touch the similarly-worded sections earlier in the document.
$n = $d.Paragraphs.Count
$base = $null
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*sub-goal*") {
        $base = $i
        break
    }
}

# Layout right after $base, before any edits:
#   base+0  b) My sub-goal is to identify a pattern and or a formula.
#   base+1  (empty paragraph holding the _GoBack bookmark)
#   base+2  (empty paragraph)
#   base+3  3) Identify potential solutions
#   base+4  a) Finding a pattern ... of 50 stop on the first finger.
#   base+5  4) Evaluate each potential solution
#   base+6  a)                                   <- needs new text + reindent
#   base+7  b) Using this pattern works ...
#   base+8  5) Choose[tab]a[tab]solution...       <- needs rewritten text
#   base+9  a) Explain...full.                    <- delete
#   base+10 b) Describe...include                 <- delete
#   base+11 drawings...clearly                    <- delete
#   base+12 communicating...solution).            <- delete

# ---------------------------------------------------------------
# 1) Remove the obsolete trailing boilerplate paragraphs.
# ---------------------------------------------------------------
$pExplain = $base + 9
$pCommunicating = $base + 12
$rKill = $d.Range($d.Paragraphs($pExplain).Range.Start, $d.Paragraphs($pCommunicating).Range.End)
$rKill.Delete()

# ---------------------------------------------------------------
# 2) Rewrite the "5) Choose a solution..." paragraph text.
# ---------------------------------------------------------------
$pChoose = $base + 8
$d.Paragraphs($pChoose).Range.Delete()
$d.Paragraphs($pChoose - 1).Range.InsertParagraphAfter()
$d.Paragraphs($pChoose).Range.InsertAfter("5) Choose a solution and develop a")
$d.Paragraphs($pChoose).Range.InsertAfter([char]9)
$d.Paragraphs($pChoose).Range.InsertAfter("plan to")
$d.Paragraphs($pChoose).Range.InsertAfter([char]9)
$d.Paragraphs($pChoose).Range.InsertAfter("implement it.")

# ---------------------------------------------------------------
# 3) Insert the two new step-5 sub-paragraphs right after it.
# ---------------------------------------------------------------
$d.Paragraphs($pChoose).Range.InsertParagraphAfter()
$pPlanA = $pChoose + 1
$d.Paragraphs($pPlanA).Format.LeftIndent = 36
$d.Paragraphs($pPlanA).Range.InsertAfter("a) Using the pattern I identified in 3a, I have determined that counts of the first 10 stop on the first finger and counts in multitudes of 50 and 100 stop on the first finger. This will tell us that counts of 10, 100 and 1000 will all stop on the first finger.")

$d.Paragraphs($pPlanA).Range.InsertParagraphAfter()
$pPlanB = $pPlanA + 1
$d.Paragraphs($pPlanB).Format.LeftIndent = 36
$d.Paragraphs($pPlanB).Range.InsertAfter("b) My hand was the best test tool for this problem. I counted out ")
$d.Paragraphs($pPlanB).Range.InsertAfter("10, then 20, then 30, then 40, then 50 and so on up to 100. I took note of what finger I stopped on at each plateau. This enabled me to verify a consistent repetitive pattern.")

# ---------------------------------------------------------------
# 4) Complete step 4a) - add the missing evaluation sentence and
#    fix the indentation (firstLine -> left, matching 3a/3b/5a/5b).
# ---------------------------------------------------------------
$pEval4a = $base + 6
$d.Paragraphs($pEval4a).Range.Delete()
$d.Paragraphs($pEval4a - 1).Range.InsertParagraphAfter()
$d.Paragraphs($pEval4a).Format.LeftIndent = 36
$d.Paragraphs($pEval4a).Range.InsertAfter("a) ")
$d.Paragraphs($pEval4a).Range.InsertAfter("While developing a formula is a viable solution, identifying a pattern will meet the goals for this problem.")

# ---------------------------------------------------------------
# 5) The pair of paragraphs right after "My sub-goal..." (base+1,
#    base+2): base+1 currently holds the "_GoBack" bookmark and
#    needs to become a plain empty paragraph; the bookmark itself
#    moves to the very end of the document (after the new 5b
#    paragraph written above).
# ---------------------------------------------------------------
$pBookmark = $base + 1
$pEmpty = $base + 2
$rTwo = $d.Range($d.Paragraphs($pBookmark).Range.Start, $d.Paragraphs($pEmpty).Range.End)
$rTwo.Delete()
$prev = $d.Paragraphs($pBookmark - 1)
$prev.Range.InsertParagraphAfter()
$prev.Range.InsertParagraphAfter()

# ---------------------------------------------------------------
# 6) Append the (now free-standing) "_GoBack" bookmark paragraph
#    to the very end of the document.
# ---------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$d.Bookmarks.Add("_GoBack", $d.Paragraphs($d.Paragraphs.Count).Range)
